$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New parts rows (34-36)
$ws.Range("A34").Value = "Molex Micro-Fit 4 Circuit plug Housing"
$ws.Range("E34").Value = "538-43020-0400"

$ws.Range("A35").Value = "Molex Micro-Fit Pins"
$ws.Range("E35").Value = "538-43031-0007"
$ws.Range("E35").WrapText = $true

$ws.Range("A36").Value = "Battery 14500 3.7V Li-Ion"
$ws.Range("F36").Value = "http://www.dx.com/p/trustfire-protected-14500-3-7v-900mah-lithium-batteries-2-pack-blue-19626"

# Match style of column A part-name cells (style s="3" used by A2:A31)
$ws.Range("A34").Style = $ws.Range("A30").Style
$ws.Range("A35").Style = $ws.Range("A30").Style
$ws.Range("A36").Style = $ws.Range("A30").Style

# Add hyperlink on F36 (becomes rId5)
$ws.Hyperlinks.Add($ws.Range("F36"), "http://www.dx.com/p/trustfire-protected-14500-3-7v-900mah-lithium-batteries-2-pack-blue-19626")

# Column A width widened to fit new longer text
$ws.Columns.Item(1).ColumnWidth = 35.28515625

# Selection / view update
$ws.Range("F36").Select()
